$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The oldest reporting period (6-month period ending 1399/06, originally column D)
# is retired from the cumulative table, and a brand-new period (12-month period
# ending 1401/12) is appended as the new last column. Deleting column D shifts
# every later column one step to the left, which is exactly what's needed here.
$ws.Columns("D").Delete()

# Clone column L's formatting (styles + column width) into the new, now-empty
# column M so the appended period matches the look of the rest of the table.
$ws.Range("L1:L28").Copy()
$ws.Range("M1:M28").PasteSpecial(-4122)
$ws.Columns("M").ColumnWidth = 28.17

# Header row: the new column's financial period label
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"

# Publish-date row: the new column's publish date, plus a correction to the
# date that slid into I9 (was "1401-10-28 (7)", corrected to "1402-02-27 (9)")
$ws.Range("I9").Value = "1402-02-27 (9)"
$ws.Range("M9").Value = "1402-02-27 (2)"

# Blank-marker rows use the literal "-" placeholder like their neighbours
$ws.Range("M15").Value = "-"
$ws.Range("M21").Value = "-"
$ws.Range("M23").Value = "-"

# New financial figures for the appended 12-month-ending-1401/12 column
$ws.Range("M11").Value = 224666
$ws.Range("M12").Value = -72185
$ws.Range("M13").Value = 152481
$ws.Range("M14").Value = -10670
$ws.Range("M16").Value = 9823
$ws.Range("M17").Value = 151634
$ws.Range("M18").Value = -472
$ws.Range("M19").Value = -19161
$ws.Range("M20").Value = 132001
$ws.Range("M22").Value = 132001
$ws.Range("M24").Value = 132001
$ws.Range("M25").Value = 0
$ws.Range("M26").Value = 10001
$ws.Range("M27").Value = 0
